# Updates cryptos list price/volume figures (GitHub Actions scheduled refresh).
# Cells hold the numbers as literal text (coinranking.com scrape format), so
# purely-numeric-looking values are written with a leading apostrophe to keep
# Excel from re-typing them as actual numbers (which would drop trailing zeros
# and introduce floating point noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.623.73"
$ws.Range("E2").Value = "  -1.76%  "
$ws.Range("D3").Value = "1.588.12"
$ws.Range("D5").Value = "'210.99"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -2.45%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("D10").Value = "'19.59"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("E11").Value = "  -1.67%  "
$ws.Range("D12").Value = "1.811.16"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "1.590.16"
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  -4.05%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "26.613.45"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'208.20"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("D21").Value = "'6.72"
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "'2.34"
$ws.Range("E23").Value = "  -3.69%  "
$ws.Range("D24").Value = "'8.85"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").Value = "'146.81"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'7.23"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  -3.40%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("E32").Value = "  -4.39%  "
$ws.Range("D33").Value = "'0.666"
$ws.Range("E33").Value = "  +20.57%  "
$ws.Range("E34").Value = "  -2.88%  "
$ws.Range("D35").Value = "1.306.93"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("E37").Value = "  -5.63%  "
$ws.Range("E38").Value = "  -3.26%  "
$ws.Range("E39").Value = "  -3.48%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("D41").Value = "'0.793"
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("D42").Value = "'5.37"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("E43").Value = "  -3.29%  "
$ws.Range("D44").Value = "'62.62"
$ws.Range("E44").Value = "  -4.53%  "
$ws.Range("D45").Value = "1.724.34"
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("D46").Value = "'89.57"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("D48").Value = "'0.839"
$ws.Range("E48").Value = "  -2.14%  "
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("E51").Value = "  -1.34%  "
